$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("A1").Value = "MYEAR"
$ws.Range("B1").Value = "Seleccion"
$ws.Range("C1").Value = "COACH"

# Alemania -> Alemania Federal (specific rows only, to match diff exactly)
$rowsAlemania = @(21, 34, 63, 79, 159)
foreach ($r in $rowsAlemania) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value2 -eq "Alemania") {
        $cell.Value = "Alemania Federal"
    }
}

# Rumanía -> Rumania (specific rows only)
$rowsRumania = @(25, 41)
foreach ($r in $rowsRumania) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value2 -eq "Rumanía") {
        $cell.Value = "Rumania"
    }
}
